$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.814.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.354.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.344.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.14%  "

# Row 9
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("E10").Value = "  +6.77%  "

# Row 11
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000281"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.884.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.121"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.343.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.729.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.12%  "

# Row 31
$ws.Range("E31").Value = "  -2.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "63.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.05%  "

# Row 33
$ws.Range("E33").Value = "  -0.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "576.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.11%  "

# Row 35
$ws.Range("E35").Value = "  -0.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.73%  "

# Row 38
$ws.Range("E38").Value = "  -0.42%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.370"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.81%  "

# Row 41
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0738"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.071.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0416"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.72%  "

# Row 44
$ws.Range("E44").Value = "  -2.58%  "

# Row 45
$ws.Range("E45").Value = "  +3.51%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.56%  "

# Row 50
$ws.Range("E50").Value = "  -2.67%  "

# Row 51
$ws.Range("E51").Value = "  -0.12%  "
